$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 644.6429000000001
$ws.Range("I28").Value = 604.3333
$ws.Range("J28").Value = 674.875
$ws.Range("K28").Value = 604.3333
$ws.Range("L28").Value = 674.875
$ws.Range("M28").Value = -119.3333
$ws.Range("N28").Value = -1644.875

$ws.Range("H51").Value = 15288.611
$ws.Range("I51").Value = 2831.6667
$ws.Range("K51").Value = 2831.6667
$ws.Range("M51").Value = -2347.6667

$ws.Range("H62").Value = 26669916
$ws.Range("I62").Value = 44445780
$ws.Range("J62").Value = 6124.75
$ws.Range("K62").Value = 44445780
$ws.Range("L62").Value = 6124.75
$ws.Range("M62").Value = -44445156
$ws.Range("N62").Value = -7372.75

$ws.Range("H64").Value = 10737.875
$ws.Range("J64").Value = 8749.25
$ws.Range("L64").Value = 8749.25
$ws.Range("N64").Value = -9245.25

$ws.Range("H65").Value = 26669916
$ws.Range("I65").Value = 44445780
$ws.Range("J65").Value = 6124.75
$ws.Range("K65").Value = 222228900
$ws.Range("L65").Value = 30623.75
$ws.Range("M65").Value = -222225780
$ws.Range("N65").Value = -36863.75

$ws.Range("H67").Value = 10737.875
$ws.Range("J67").Value = 8749.25
$ws.Range("L67").Value = 8749.25
$ws.Range("N67").Value = -10465.25

$ws.Range("H88").Value = 6195.3125
$ws.Range("J88").Value = 6344.1665
$ws.Range("L88").Value = 6344.1665
$ws.Range("N88").Value = -7156.1665

$ws.Range("H91").Value = 6195.3125
$ws.Range("J91").Value = 6344.1665
$ws.Range("L91").Value = 6344.1665
$ws.Range("N91").Value = -9152.166499999999

$ws.Range("H132").Value = 1824.1428
$ws.Range("I132").Value = 1824.1428
$ws.Range("K132").Value = 5472.428400000001
$ws.Range("M132").Value = -2942.428400000001

$ws.Range("H137").Value = 22736710
$ws.Range("I137").Value = 50001600
$ws.Range("J137").Value = 15967.167
$ws.Range("K137").Value = 150004800
$ws.Range("L137").Value = 47901.501
$ws.Range("M137").Value = -150002250
$ws.Range("N137").Value = -53001.501

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 19500
$ws.Range("I55").Value = 8500
$ws.Range("K55").Value = 8500
$ws.Range("M55").Value = -8185

$ws.Range("H61").Value = 5563.25
$ws.Range("I61").Value = 4064.9443
$ws.Range("J61").Value = 10058.167
$ws.Range("K61").Value = 4064.9443
$ws.Range("L61").Value = 10058.167
$ws.Range("M61").Value = -3852.9443
$ws.Range("N61").Value = -10482.167

$ws.Range("H63").Value = 8157.933
$ws.Range("J63").Value = 9425
$ws.Range("L63").Value = 9425
$ws.Range("N63").Value = -10797

$ws.Range("H66").Value = 8157.933
$ws.Range("J66").Value = 9425
$ws.Range("L66").Value = 47125
$ws.Range("N66").Value = -53989

$ws.Range("H94").Value = 46332.668
$ws.Range("J94").Value = 46332.668
$ws.Range("L94").Value = 46332.668
$ws.Range("N94").Value = -48134.668

$ws.Range("H136").Value = 5563.25
$ws.Range("I136").Value = 4064.9443
$ws.Range("J136").Value = 10058.167
$ws.Range("K136").Value = 12194.8329
$ws.Range("L136").Value = 30174.501
$ws.Range("M136").Value = -9644.832900000001
$ws.Range("N136").Value = -35274.501

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 56380.945
$ws.Range("I80").Value = 574.25
$ws.Range("J80").Value = 101026.3
$ws.Range("K80").Value = 574.25
$ws.Range("L80").Value = 101026.3
$ws.Range("M80").Value = 423.75
$ws.Range("N80").Value = -103022.3

$ws.Range("H83").Value = 56380.945
$ws.Range("I83").Value = 574.25
$ws.Range("J83").Value = 101026.3
$ws.Range("K83").Value = 2871.25
$ws.Range("L83").Value = 505131.5
$ws.Range("M83").Value = 2120.75
$ws.Range("N83").Value = -515115.5

$ws.Range("H105").Value = 50002000
$ws.Range("I105").Value = 62501892
$ws.Range("K105").Value = 62501892
$ws.Range("M105").Value = -62500145

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H103").Value = 23498.25
$ws.Range("J103").Value = 79980
$ws.Range("L103").Value = 79980
$ws.Range("N103").Value = -82324

$ws.Range("H104").Value = 56666.332
$ws.Range("J104").Value = 56666.332
$ws.Range("L104").Value = 56666.332
$ws.Range("N104").Value = -61908.332

$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H20").Value = 3098.5
$ws.Range("J20").Value = 3098.5
$ws.Range("L20").Value = 9295.5
$ws.Range("N20").Value = -9749.5

$ws.Range("H80").Value = 3825.5
$ws.Range("I80").Value = 3002
$ws.Range("K80").Value = 9006
$ws.Range("M80").Value = -8070

$ws.Range("H83").Value = 3825.5
$ws.Range("I83").Value = 3002
$ws.Range("K83").Value = 27018
$ws.Range("M83").Value = -22338

$ws.Range("H107").Value = 1201.1471
$ws.Range("J107").Value = 1357.8966
$ws.Range("L107").Value = 4073.6898
$ws.Range("N107").Value = -7913.6898

$ws.Range("H122").Value = 76929970
$ws.Range("I122").Value = 111119960
$ws.Range("K122").Value = 1000079640
$ws.Range("M122").Value = -1000077190

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 44998.5
$ws.Range("J93").Value = 44998.5
$ws.Range("L93").Value = 44998.5
$ws.Range("N93").Value = -48742.5

$ws.Range("H107").Value = 312.2
$ws.Range("I107").Value = 325.66666
$ws.Range("K107").Value = 325.66666
$ws.Range("M107").Value = 1594.33334

$ws.Range("H113").Value = 40963.89
$ws.Range("I113").Value = 3113.3333
$ws.Range("J113").Value = 116665
$ws.Range("K113").Value = 3113.3333
$ws.Range("L113").Value = 116665
$ws.Range("M113").Value = -943.3332999999998
$ws.Range("N113").Value = -121005

$ws.Range("H122").Value = 6816.963
$ws.Range("I122").Value = 7432.3687
$ws.Range("K122").Value = 22297.1061
$ws.Range("M122").Value = -19847.1061

$ws.Range("H126").Value = 3175.9546
$ws.Range("I126").Value = 2579.9285
$ws.Range("J126").Value = 4219
$ws.Range("K126").Value = 7739.7855
$ws.Range("L126").Value = 12657
$ws.Range("M126").Value = -5269.7855
$ws.Range("N126").Value = -17597

$ws.Range("H132").Value = 9192.277
$ws.Range("I132").Value = 6033.273
$ws.Range("K132").Value = 18099.819
$ws.Range("M132").Value = -15569.819

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").ClearContents()

$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("M28").ClearContents()
$ws.Range("N28").ClearContents()

$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("M37").ClearContents()
$ws.Range("N37").ClearContents()

$ws.Range("H46").Value = 8169.3076
$ws.Range("I46").Value = 3199.6667
$ws.Range("J46").Value = 9660.200000000001
$ws.Range("K46").Value = 3199.6667
$ws.Range("L46").Value = 9660.200000000001
$ws.Range("M46").Value = -3011.6667
$ws.Range("N46").Value = -10036.2

$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 71999.75
$ws.Range("J43").Value = 89333
$ws.Range("L43").Value = 89333
$ws.Range("N43").Value = -89631

$ws.Range("H54").Value = 42484.625
$ws.Range("I54").Value = 79000
$ws.Range("J54").Value = 20575.4
$ws.Range("K54").Value = 79000
$ws.Range("L54").Value = 20575.4
$ws.Range("M54").Value = -78480
$ws.Range("N54").Value = -21615.4

$ws.Range("H132").Value = 5129.8823
$ws.Range("I132").Value = 1957.5625
$ws.Range("K132").Value = 5872.6875
$ws.Range("M132").Value = -3342.6875
